$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 08:57:42'
$ws.Cells.Item(3, 1).Value = 'Total filas: 85'
$ws.Cells.Item(15, 3).Value = '215A_EL PATO'
$ws.Cells.Item(16, 3).Value = '225_GOMEZ'
$ws.Cells.Item(55, 1).Value = '08:49:06'
$ws.Cells.Item(55, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(55, 4).Value = 4
$ws.Cells.Item(56, 1).Value = '08:14:55'
$ws.Cells.Item(56, 3).Value = '215B_EL PATO'
$ws.Cells.Item(56, 4).Value = 39
$ws.Cells.Item(59, 1).Value = '08:57:42'
$ws.Cells.Item(59, 4).Value = 1
$ws.Cells.Item(61, 1).Value = '08:57:42'
$ws.Cells.Item(61, 4).Value = 8
$ws.Cells.Item(62, 1).Value = '08:57:42'
$ws.Cells.Item(62, 4).Value = 9
$ws.Cells.Item(63, 1).Value = '08:57:42'
$ws.Cells.Item(63, 4).Value = 19
$ws.Cells.Item(64, 1).Value = '08:57:42'
$ws.Cells.Item(64, 4).Value = 20
$ws.Cells.Item(66, 1).Value = '08:57:42'
$ws.Cells.Item(66, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(66, 4).Value = 21
$ws.Cells.Item(67, 3).Value = '14_ABASTO'
$ws.Cells.Item(71, 1).Value = '08:14:55'
$ws.Cells.Item(71, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(71, 4).Value = 77
$ws.Cells.Item(72, 1).Value = '08:49:06'
$ws.Cells.Item(72, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(72, 4).Value = 42
$ws.Cells.Item(75, 1).Value = '08:57:42'
$ws.Cells.Item(75, 2).Value = '09:36'
$ws.Cells.Item(75, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(75, 4).Value = 39
$ws.Cells.Item(76, 1).Value = '08:57:42'
$ws.Cells.Item(76, 2).Value = '09:39'
$ws.Cells.Item(76, 3).Value = '15_ABASTO'
$ws.Cells.Item(76, 4).Value = 42
$ws.Cells.Item(77, 1).Value = '08:57:42'
$ws.Cells.Item(77, 2).Value = '09:41'
$ws.Cells.Item(77, 4).Value = 44
$ws.Cells.Item(78, 2).Value = '09:42'
$ws.Cells.Item(78, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(78, 4).Value = 53
$ws.Cells.Item(79, 1).Value = '08:57:42'
$ws.Cells.Item(79, 2).Value = '09:43'
$ws.Cells.Item(79, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(79, 4).Value = 46
$ws.Cells.Item(80, 1).Value = '08:57:42'
$ws.Cells.Item(80, 2).Value = '09:53'
$ws.Cells.Item(80, 3).Value = '10_OLMOS'
$ws.Cells.Item(80, 4).Value = 56
$ws.Cells.Item(81, 1).Value = '08:57:42'
$ws.Cells.Item(81, 2).Value = '09:58'
$ws.Cells.Item(81, 4).Value = 61
$ws.Cells.Item(82, 1).Value = '08:49:06'
$ws.Cells.Item(82, 2).Value = '09:59'
$ws.Cells.Item(82, 3).Value = '215C_EL PATO'
$ws.Cells.Item(82, 4).Value = 70
$ws.Cells.Item(83, 1).Value = '08:57:42'
$ws.Cells.Item(83, 2).Value = '10:05'
$ws.Cells.Item(83, 4).Value = 68
$ws.Cells.Item(84, 2).Value = '10:06'
$ws.Cells.Item(84, 3).Value = '14_ABASTO'
$ws.Cells.Item(84, 4).Value = 77
$ws.Cells.Item(85, 1).Value = '08:57:42'
$ws.Cells.Item(85, 2).Value = '10:13'
$ws.Cells.Item(85, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(85, 4).Value = 76
$ws.Cells.Item(86, 1).Value = '08:57:42'
$ws.Cells.Item(86, 2).Value = '10:24'
$ws.Cells.Item(86, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(86, 4).Value = 87
$ws.Cells.Item(87, 1).Value = '08:57:42'
$ws.Cells.Item(87, 2).Value = '10:25'
$ws.Cells.Item(87, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(87, 4).Value = 88
$ws.Cells.Item(88, 1).Value = '08:57:42'
$ws.Cells.Item(88, 2).Value = '10:29'
$ws.Cells.Item(88, 3).Value = '15_ABASTO'
$ws.Cells.Item(88, 4).Value = 92
$ws.Cells.Item(89, 1).Value = '08:57:42'
$ws.Cells.Item(89, 2).Value = '10:44'
$ws.Cells.Item(89, 3).Value = '11X44_ETCHEVERRY'
$ws.Cells.Item(89, 4).Value = 107
$ws.Cells.Item(89, 5).Value = 'LP1912'
$ws.Cells.Item(90, 1).Value = '08:57:42'
$ws.Cells.Item(90, 2).Value = '10:46'
$ws.Cells.Item(90, 3).Value = '15_P INDUSTRIAL'
$ws.Cells.Item(90, 4).Value = 109
$ws.Cells.Item(90, 5).Value = 'LP1912'

$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 08:57:42'
$ws.Cells.Item(19, 1).Value = '08:57:42'
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(20, 1).Value = '08:57:42'
$ws.Cells.Item(20, 4).Value = 61

$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 08:57:42'
$ws.Cells.Item(3, 1).Value = 'Total filas: 16'
$ws.Cells.Item(13, 1).Value = '08:57:42'
$ws.Cells.Item(13, 2).Value = '08:59'
$ws.Cells.Item(13, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(13, 4).Value = 2
$ws.Cells.Item(13, 5).Value = 'L6203'
$ws.Cells.Item(14, 1).Value = '08:57:42'
$ws.Cells.Item(14, 2).Value = '09:20'
$ws.Cells.Item(14, 4).Value = 23
$ws.Cells.Item(15, 1).Value = '08:49:06'
$ws.Cells.Item(15, 2).Value = '09:21'
$ws.Cells.Item(15, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(15, 4).Value = 32
$ws.Cells.Item(15, 5).Value = 'L6173'
$ws.Cells.Item(16, 1).Value = '08:57:42'
$ws.Cells.Item(16, 2).Value = '10:12'
$ws.Cells.Item(16, 4).Value = 75
$ws.Cells.Item(17, 2).Value = '10:13'
$ws.Cells.Item(17, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(17, 4).Value = 84
$ws.Cells.Item(17, 5).Value = 'L6203'
$ws.Cells.Item(18, 1).Value = '08:57:42'
$ws.Cells.Item(18, 2).Value = '10:29'
$ws.Cells.Item(18, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(18, 4).Value = 92
$ws.Cells.Item(19, 1).Value = '08:57:42'
$ws.Cells.Item(19, 2).Value = '10:30'
$ws.Cells.Item(19, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(19, 4).Value = 93
$ws.Cells.Item(19, 5).Value = 'L6173'
$ws.Cells.Item(20, 1).Value = '08:49:06'
$ws.Cells.Item(20, 2).Value = '10:30'
$ws.Cells.Item(20, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(20, 4).Value = 101
$ws.Cells.Item(20, 5).Value = 'L6173'
$ws.Cells.Item(21, 1).Value = '08:49:06'
$ws.Cells.Item(21, 2).Value = '10:31'
$ws.Cells.Item(21, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(21, 4).Value = 102
$ws.Cells.Item(21, 5).Value = 'L6173'
